$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.042.66'
$ws.Range("E2").Value = '  -2.70%  '

$ws.Range("D3").Value = '''3.206.70'
$ws.Range("E3").Value = '  -1.17%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").Value = '''602.11'
$ws.Range("E5").Value = '  -0.58%  '

$ws.Range("D6").Value = '''154.03'
$ws.Range("E6").Value = '  -2.45%  '

$ws.Range("E7").Value = '  +0.15%  '

$ws.Range("D8").Value = '''3.205.61'
$ws.Range("E8").Value = '  -1.14%  '

$ws.Range("D9").Value = '''0.530'
$ws.Range("E9").Value = '  -3.38%  '

$ws.Range("D10").Value = '''0.155'
$ws.Range("E10").Value = '  -3.60%  '

$ws.Range("D11").Value = '''5.56'
$ws.Range("E11").Value = '  -2.14%  '

$ws.Range("E12").Value = '  -5.77%  '

$ws.Range("D13").Value = '''0.0000259'
$ws.Range("E13").Value = '  -4.83%  '

$ws.Range("D14").Value = '''37.49'
$ws.Range("E14").Value = '  -4.22%  '

$ws.Range("D15").Value = '''3.736.93'
$ws.Range("E15").Value = '  -1.05%  '

$ws.Range("D16").Value = '''65.209.30'

$ws.Range("D17").Value = '''3.208.94'
$ws.Range("E17").Value = '  -1.00%  '

$ws.Range("E18").Value = '  +0.48%  '

$ws.Range("D19").Value = '''7.11'
$ws.Range("E19").Value = '  -5.71%  '

$ws.Range("D20").Value = '''486.84'
$ws.Range("E20").Value = '  -5.23%  '

$ws.Range("D21").Value = '''14.98'
$ws.Range("E21").Value = '  -2.79%  '

$ws.Range("E22").Value = '  -1.64%  '

$ws.Range("D23").Value = '''7.82'
$ws.Range("E23").Value = '  -3.65%  '

$ws.Range("D24").Value = '''14.05'
$ws.Range("E24").Value = '  -5.97%  '

$ws.Range("D25").Value = '''85.05'
$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").Value = '''0.996'
$ws.Range("E26").Value = '  -0.52%  '

$ws.Range("E27").Value = '  -1.36%  '

$ws.Range("D28").Value = '''8.81'
$ws.Range("E28").Value = '  -5.74%  '

$ws.Range("E29").Value = '  +38.08%  '

$ws.Range("E30").Value = '  -5.00%  '

$ws.Range("D31").Value = '''6.93'
$ws.Range("E31").Value = '  -2.38%  '

$ws.Range("E32").Value = '  -8.76%  '

$ws.Range("D33").Value = '''27.22'
$ws.Range("E33").Value = '  -4.11%  '

$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  +0.10%  '

$ws.Range("D35").Value = '''1.11'
$ws.Range("E35").Value = '  -6.02%  '

$ws.Range("D36").Value = '''6.19'
$ws.Range("E36").Value = '  -5.75%  '

$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").Value = '''3.31'
$ws.Range("E37").Value = '  +9.43%  '

$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = '''54.79'
$ws.Range("E38").Value = '  -2.69%  '

$ws.Range("D39").Value = '''481.12'
$ws.Range("E39").Value = '  -8.01%  '

$ws.Range("D40").Value = '''0.0₃0745'
$ws.Range("E40").Value = '  -2.85%  '

$ws.Range("D41").Value = '''0.0409'
$ws.Range("E41").Value = '  -2.94%  '

$ws.Range("D42").Value = '''0.126'
$ws.Range("E42").Value = '  -1.91%  '

$ws.Range("E43").Value = '  -2.98%  '

$ws.Range("D44").Value = '''2.48'
$ws.Range("E44").Value = '  -0.83%  '

$ws.Range("D45").Value = '''2.946.05'
$ws.Range("E45").Value = '  +2.64%  '

$ws.Range("D46").Value = '''0.282'
$ws.Range("E46").Value = '  -6.93%  '

$ws.Range("D47").Value = '''27.68'
$ws.Range("E47").Value = '  -3.26%  '

$ws.Range("E48").Value = '  -1.87%  '

$ws.Range("E49").Value = '  -0.03%  '

$ws.Range("E50").Value = '  -0.42%  '

$ws.Range("D51").Value = '''121.10'
